$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# New matchup rows appended to the Nine-ball data table (rows 860-876, cols A:D)
$data = @(
    @(5,18,6,2),
    @(4,5,5,15),
    @(4,5,5,15),
    @(4,14,3,6),
    @(4,13,3,7),
    @(5,12,4,8),
    @(6,12,5,8),
    @(4,14,3,6),
    @(4,18,5,2),
    @(3,12,4,8),
    @(5,15,7,5),
    @(7,6,5,14),
    @(3,8,4,12),
    @(4,13,2,7),
    @(2,12,3,8),
    @(2,18,3,2),
    @(9,7,6,13)
)

$startRow = 860
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

# Scroll/select to reflect the new bottom of the sheet, like a user who just
# pasted the new rows and left the selection just past the last one.
$excel.ActiveWindow.ScrollRow = 866
$ws.Range("F878").Select()
